$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 36; this shifts rows 36-77 down to 37-78 and
# Excel automatically re-points every formula reference (B39->B40, L69->L70, ...)
$ws.Rows("36:36").Insert()

# Fill the new row's content
$ws.Range("A36").Value = "solver"
$ws.Range("B36").Value = "gurobi"

# Style B36: numeric format id 164 ("yyyy-mm-dd hh:mm:ss") + right alignment
$ws.Range("B36").NumberFormat = "yyyy\-mm\-dd\ hh:mm:ss"
$ws.Range("B36").HorizontalAlignment = -4152

# Style C36: numeric format id 164 only (no alignment change), empty cell
$ws.Range("C36").NumberFormat = "yyyy\-mm\-dd\ hh:mm:ss"

# Update the selection to match the post-edit state
$ws.Range("C36").Select()
